# Updates the cryptos price/volume table with the latest scraped values.
# Note: values in columns D (Price) and E (Volume(1h)) are text, not numbers
# (prices use "." as both thousands and decimal separators, e.g. "26.733.16",
# and volumes are padded percentage strings like "  +3.77%  "). A leading
# apostrophe forces Excel to store them as text instead of auto-converting
# ambiguous-looking values (e.g. "1.001") into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.733.16'
$ws.Range("E2").Value = '''  +3.77%  '

$ws.Range("D3").Value = '''1.868.26'
$ws.Range("E3").Value = '''  +2.82%  '

$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '''  +0.10%  '

$ws.Range("D5").Value = '''277.05'
$ws.Range("E5").Value = '''  -0.66%  '

$ws.Range("D6").Value = '''1.001'
$ws.Range("E6").Value = '''  +0.12%  '

$ws.Range("D7").Value = '''0.5296'
$ws.Range("E7").Value = '''  +3.92%  '

$ws.Range("D8").Value = '''0.3419'
$ws.Range("E8").Value = '''  -3.50%  '

$ws.Range("D9").Value = '''0.06918'
$ws.Range("E9").Value = '''  +3.75%  '

$ws.Range("D10").Value = '''20.02'

$ws.Range("D11").Value = '''0.8025'
$ws.Range("E11").Value = '''  -2.97%  '

$ws.Range("D12").Value = '''0.07739'

$ws.Range("D13").Value = '''1.874.07'
$ws.Range("E13").Value = '''  +3.47%  '

$ws.Range("D14").Value = '''90.12'
$ws.Range("E14").Value = '''  +2.65%  '

$ws.Range("D15").Value = '''5.166'
$ws.Range("E15").Value = '''  +1.75%  '

$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").Value = '''14.55'
$ws.Range("E16").Value = '''  +3.19%  '

$ws.Range("B17").Value = 'BinanceUSD'
$ws.Range("C17").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D17").Value = '''0.9996'
$ws.Range("E17").Value = '''  +0.00%  '

$ws.Range("E18").Value = '''  -0.08%  '

$ws.Range("D20").Value = '''26.769.31'
$ws.Range("E20").Value = '''  +3.76%  '

$ws.Range("D21").Value = '''2.120.59'
$ws.Range("E21").Value = '''  +4.66%  '

$ws.Range("D22").Value = '''4.746'
$ws.Range("E22").Value = '''  -0.12%  '

$ws.Range("D23").Value = '''10.02'
$ws.Range("E23").Value = '''  +0.23%  '

$ws.Range("D24").Value = '''6.177'
$ws.Range("E24").Value = '''  +0.92%  '

$ws.Range("D25").Value = '''2.363'
$ws.Range("E25").Value = '''  +5.88%  '

$ws.Range("D26").Value = '''146.07'
$ws.Range("E26").Value = '''  +2.63%  '

$ws.Range("D27").Value = '''17.31'
$ws.Range("E27").Value = '''  +0.99%  '

$ws.Range("D28").Value = '''1.653'
$ws.Range("E28").Value = '''  -0.85%  '

$ws.Range("D29").Value = '''113.03'
$ws.Range("E29").Value = '''  +3.35%  '

$ws.Range("D30").Value = '''4.328'
$ws.Range("E30").Value = '''  -0.23%  '

$ws.Range("D31").Value = '''4.335'
$ws.Range("E31").Value = '''  +2.26%  '

$ws.Range("D32").Value = '''0.08894'
$ws.Range("E32").Value = '''  +1.39%  '

$ws.Range("D33").Value = '''0.04946'
$ws.Range("E33").Value = '''  +0.87%  '

$ws.Range("D34").Value = '''1.164'
$ws.Range("E34").Value = '''  +2.12%  '

$ws.Range("D35").Value = '''0.7283'
$ws.Range("E35").Value = '''  -0.56%  '

$ws.Range("D36").Value = '''2.880'
$ws.Range("E36").Value = '''  +0.31%  '

$ws.Range("D37").Value = '''3.251'
$ws.Range("E37").Value = '''  +3.49%  '

$ws.Range("D38").Value = '''0.01855'
$ws.Range("E38").Value = '''  +0.10%  '

$ws.Range("D39").Value = '''2.322'
$ws.Range("E39").Value = '''  -3.03%  '

$ws.Range("D40").Value = '''0.5144'
$ws.Range("E40").Value = '''  -0.29%  '

$ws.Range("D41").Value = '''0.9470'
$ws.Range("E41").Value = '''  -1.92%  '

$ws.Range("D42").Value = '''116.44'
$ws.Range("E42").Value = '''  +4.80%  '

$ws.Range("D43").Value = '''6.141'
$ws.Range("E43").Value = '''  -1.32%  '

$ws.Range("D44").Value = '''8.093'
$ws.Range("E44").Value = '''  +0.50%  '

$ws.Range("E45").Value = '''  +0.09%  '

$ws.Range("D46").Value = '''0.4442'
$ws.Range("E46").Value = '''  -2.59%  '

$ws.Range("D47").Value = '''0.1339'
$ws.Range("E47").Value = '''  -2.21%  '

$ws.Range("D48").Value = '''9.330'
$ws.Range("E48").Value = '''  +1.48%  '

$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").Value = '''36.34'
$ws.Range("E49").Value = '''  -0.70%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '''0.06017'
$ws.Range("E50").Value = '''  +3.25%  '

$ws.Range("D51").Value = '''1.874.07'
$ws.Range("E51").Value = '''  +3.47%  '
